# Update the F1_Venta_23_Ene_Porcentaje results (column C) on Hoja1 with the
# final percentages, replacing the placeholder/interim values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = 0.47
$ws.Range("C3").Value  = 0.46
$ws.Range("C4").Value  = 0.56
$ws.Range("C5").Value  = 0.2
$ws.Range("C6").Value  = 0.3
$ws.Range("C7").Value  = 0.48
$ws.Range("C8").Value  = 0.31
$ws.Range("C9").Value  = 0.44
$ws.Range("C10").Value = 0.43
$ws.Range("C11").Value = 0.36
$ws.Range("C12").Value = 0.44
$ws.Range("C13").Value = 0.24

# Leave the active cell/selection on C11, matching where the user last clicked.
$ws.Range("C11").Select()
